$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.731256310823452
$ws.Range("C2").Value = 5.285850051013121
$ws.Range("E2").Value = 16.48685598771166
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 3.620103850911608
$ws.Range("I2").Value = 19.05331118567647
$ws.Range("K2").Value = 8.493242427114655
$ws.Range("N2").Value = 17.47350436778733
$ws.Range("O2").Value = 20.01887315098602
$ws.Range("B3").Value = 8.402703154092155
$ws.Range("C3").Value = 5.064474738559142
$ws.Range("E3").Value = 15.55380046419533
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 3.621776292545836
$ws.Range("I3").Value = 19.14799814755978
$ws.Range("K3").Value = 8.266680856509749
$ws.Range("N3").Value = 17.53001518415936
$ws.Range("O3").Value = 20.10083240152243
$ws.Range("B4").Value = 8.19549823733208
$ws.Range("C4").Value = 4.922436083984493
$ws.Range("E4").Value = 14.95600349312518
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 3.622856248355717
$ws.Range("I4").Value = 19.2099974006202
$ws.Range("K4").Value = 8.125655209943549
$ws.Range("N4").Value = 17.56639366339047
$ws.Range("O4").Value = 20.15590539820283
$ws.Range("B5").Value = 8.10981548545875
$ws.Range("C5").Value = 4.863063140124837
$ws.Range("E5").Value = 14.70640538831743
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 3.623309727855603
$ws.Range("I5").Value = 19.23623274716568
$ws.Range("K5").Value = 8.067790647899823
$ws.Range("N5").Value = 17.58164205054554
$ws.Range("O5").Value = 20.17953871921995
$ws.Range("B6").Value = 8.095516824485314
$ws.Range("C6").Value = 4.853115737289784
$ws.Range("E6").Value = 14.66460689533891
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 3.623385837718741
$ws.Range("I6").Value = 19.24064768008242
$ws.Range("K6").Value = 8.058161010078328
$ws.Range("N6").Value = 17.58419967089545
$ws.Range("O6").Value = 20.18353481591756
$ws.Range("B7").Value = 8.194347548006611
$ws.Range("C7").Value = 4.921641330807557
$ws.Range("E7").Value = 14.95266118802348
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 3.622862309870342
$ws.Range("I7").Value = 19.21034729257146
$ws.Range("K7").Value = 8.124876311979856
$ws.Range("N7").Value = 17.56659759061418
$ws.Range("O7").Value = 20.15621930960555
$ws.Range("B8").Value = 8.619181223941547
$ws.Range("C8").Value = 5.210814328280755
$ws.Range("E8").Value = 16.17045195415212
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 3.620669520608454
$ws.Range("I8").Value = 19.08515742287965
$ws.Range("K8").Value = 8.415571945034738
$ws.Range("N8").Value = 17.49264124626519
$ws.Range("O8").Value = 20.04614467648594
$ws.Range("B9").Value = 9.40379011832459
$ws.Range("C9").Value = 5.727671106647951
$ws.Range("E9").Value = 18.43350206724235
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.61678854473426
$ws.Range("I9").Value = 18.87033065705826
$ws.Range("K9").Value = 8.966882288643177
$ws.Range("N9").Value = 17.36088997086348
$ws.Range("O9").Value = 19.86814058557959
$ws.Range("B10").Value = 9.944791123721302
$ws.Range("C10").Value = 6.074915331238961
$ws.Range("E10").Value = 20.06841821108013
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.61418985241171
$ws.Range("I10").Value = 18.73124296336817
$ws.Range("K10").Value = 9.356137683148024
$ws.Range("N10").Value = 17.272104285502
$ws.Range("O10").Value = 19.76066115226335
$ws.Range("B11").Value = 10.18220092699345
$ws.Range("C11").Value = 6.225531510947703
$ws.Range("E11").Value = 20.76994946461659
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.613061901632129
$ws.Range("I11").Value = 18.67205044179281
$ws.Range("K11").Value = 9.528978598031593
$ws.Range("N11").Value = 17.23343569554424
$ws.Range("O11").Value = 19.71687290455974
$ws.Range("B12").Value = 10.27078180347523
$ws.Range("C12").Value = 6.281489724320844
$ws.Range("E12").Value = 21.02956847929275
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.612642525399438
$ws.Range("I12").Value = 18.65022355671945
$ws.Range("K12").Value = 9.595987292077099
$ws.Range("N12").Value = 17.21903902066127
$ws.Range("O12").Value = 19.70102936805005
$ws.Range("B13").Value = 10.25176410932828
$ws.Range("C13").Value = 6.269486303521155
$ws.Range("E13").Value = 20.97392279794264
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.612732501348011
$ws.Range("I13").Value = 18.65489818981163
$ws.Range("K13").Value = 9.579838922213803
$ws.Range("N13").Value = 17.22212866645634
$ws.Range("O13").Value = 19.70440866461943
$ws.Range("B14").Value = 10.18951537728605
$ws.Range("C14").Value = 6.230156913234022
$ws.Range("E14").Value = 20.79142925403989
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.613027244122645
$ws.Range("I14").Value = 18.67024293705575
$ws.Range("K14").Value = 9.534321958968425
$ws.Range("N14").Value = 17.23224634296607
$ws.Range("O14").Value = 19.71555463106087
$ws.Range("B15").Value = 10.15121224945205
$ws.Range("C15").Value = 6.205925741058848
$ws.Range("E15").Value = 20.67886163639158
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.613208791210925
$ws.Range("I15").Value = 18.67971866223581
$ws.Range("K15").Value = 9.506352587732254
$ws.Range("N15").Value = 17.23847575201759
$ws.Range("O15").Value = 19.72247810266761
$ws.Range("B16").Value = 9.929094609308162
$ws.Range("C16").Value = 6.06492254693813
$ws.Range("E16").Value = 20.02172511707847
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.614264653953257
$ws.Range("I16").Value = 18.73519354025075
$ws.Range("K16").Value = 9.344751607854409
$ws.Range("N16").Value = 17.27466589289566
$ws.Range("O16").Value = 19.76362588614766
$ws.Range("B17").Value = 9.790552341909203
$ws.Range("C17").Value = 5.976524280278712
$ws.Range("E17").Value = 19.6078079839845
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.614926246308813
$ws.Range("I17").Value = 18.77027136967501
$ws.Range("K17").Value = 9.244485315640333
$ws.Range("N17").Value = 17.29730722313763
$ws.Range("O17").Value = 19.79017906759459
$ws.Range("B18").Value = 9.710051739252291
$ws.Range("C18").Value = 5.92498938070064
$ws.Range("E18").Value = 19.36575976442062
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.615311881404779
$ws.Range("I18").Value = 18.79083106824827
$ws.Range("K18").Value = 9.186419619741891
$ws.Range("N18").Value = 17.31049192163496
$ws.Range("O18").Value = 19.80593199184448
$ws.Range("B19").Value = 9.682657999937284
$ws.Range("C19").Value = 5.907422625193004
$ws.Range("E19").Value = 19.28312285189353
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.615443328799929
$ws.Range("I19").Value = 18.79785811378774
$ws.Range("K19").Value = 9.166693646401884
$ws.Range("N19").Value = 17.31498389273935
$ws.Range("O19").Value = 19.81134801741798
$ws.Range("B20").Value = 9.805385270248733
$ws.Range("C20").Value = 5.986006055283204
$ws.Range("E20").Value = 19.65228122917139
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.614855290653946
$ws.Range("I20").Value = 18.76649753652053
$ws.Range("K20").Value = 9.255200191358053
$ws.Range("N20").Value = 17.2948802561925
$ws.Range("O20").Value = 19.78730270367983
$ws.Range("B21").Value = 10.20783570983254
$ws.Range("C21").Value = 6.241738279562059
$ws.Range("E21").Value = 20.84519557779458
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.612940460899977
$ws.Range("I21").Value = 18.6657198416638
$ws.Range("K21").Value = 9.547710025585561
$ws.Range("N21").Value = 17.22926786097999
$ws.Range("O21").Value = 19.71226072590067
$ws.Range("B22").Value = 10.46313378207995
$ws.Range("C22").Value = 6.402590180586887
$ws.Range("E22").Value = 21.58966871231237
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.61173418795169
$ws.Range("I22").Value = 18.60328405587727
$ws.Range("K22").Value = 9.78077422219409
$ws.Range("N22").Value = 17.18782137159912
$ws.Range("O22").Value = 19.66752064818281
$ws.Range("B23").Value = 10.32760438051128
$ws.Range("C23").Value = 6.317321326831434
$ws.Range("E23").Value = 21.19553627638737
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.612373878027433
$ws.Range("I23").Value = 18.63629302326301
$ws.Range("K23").Value = 9.650667239658297
$ws.Range("N23").Value = 17.20981120646932
$ws.Range("O23").Value = 19.69100410431828
$ws.Range("B24").Value = 9.798681947928443
$ws.Range("C24").Value = 5.981721565863641
$ws.Range("E24").Value = 19.63218759175994
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.614887353264843
$ws.Range("I24").Value = 18.76820246241443
$ws.Range("K24").Value = 9.250357303494646
$ws.Range("N24").Value = 17.29597696484632
$ws.Range("O24").Value = 19.78860159061541
$ws.Range("B25").Value = 9.197379329582342
$ws.Range("C25").Value = 5.593437153786605
$ws.Range("E25").Value = 17.79378962551509
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.617793878957055
$ws.Range("I25").Value = 18.92515881192496
$ws.Range("K25").Value = 8.820217068425665
$ws.Range("N25").Value = 17.39511932353742
$ws.Range("O25").Value = 19.91221976678664
